# Apply cryptos list update (prices + 1h volume %) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.512.87"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.562.95"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -1.49%  "
$ws.Range("D5").Value = "'210.39"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.989"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("D8").Value = "'22.48"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "1.787.51"
$ws.Range("D13").Value = "1.568.34"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  +0.64%  "
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "27.512.38"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").Value = "'62.42"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "'224.47"
$ws.Range("E18").Value = "  +4.20%  "
$ws.Range("D19").Value = "'7.51"
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").Value = "0.0₃0706"
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'0.988"
$ws.Range("E21").Value = "  -1.55%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'9.41"
$ws.Range("E23").Value = "  +2.31%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'150.14"
$ws.Range("E25").Value = "  -1.69%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'15.16"
$ws.Range("E27").Value = "  +0.91%  "
$ws.Range("D28").Value = "'0.107"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "'0.990"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("E32").Value = "  +0.31%  "
$ws.Range("D33").Value = "1.461.66"
$ws.Range("E33").Value = "  +2.73%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +2.37%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'0.542"
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").Value = "'0.816"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.34"
$ws.Range("E42").Value = "  +1.19%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.989"
$ws.Range("E43").Value = "  -1.54%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'1.88"
$ws.Range("E44").Value = "  +7.83%  "
$ws.Range("D45").Value = "'0.976"
$ws.Range("E45").Value = "  -2.27%  "
$ws.Range("D46").Value = "'65.22"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("D47").Value = "1.702.21"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").Value = "'86.63"
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").Value = "'0.0523"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("E51").Value = "  -0.65%  "
